$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "2016-03-31 05:24:00"
$wsZhCn.Range("H3").Value = "2016-03-31 05:24:55"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "2016-03-31 05:24:11"
$wsDeDe.Range("H3").Value = "2016-03-31 05:25:11"
